$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark the date cell as text first so Excel doesn't auto-convert the
# "yyyy/mm/dd"-looking string into a date serial number.
$ws.Range("A43").NumberFormat = "@"

# Append the new row of mod data (row 43)
$ws.Cells.Item(43, 1).Value = "2025/12/23"
$ws.Cells.Item(43, 2).Value = "逃离鸭科夫"
$ws.Cells.Item(43, 3).Value = 1095

# Drop the temporary text-number-format override and re-apply the same
# center/center alignment used by the rest of the data rows (e.g. row 42).
$ws.Range("A43").ClearFormats()
$ws.Range("A43:C43").HorizontalAlignment = -4108
$ws.Range("A43:C43").VerticalAlignment = -4108
